# Quarterly indexing bug-fix: column A held the quarter "start" date
# (the 1st of the quarter's first month) for every forecast row, but the
# naive QoQ export is supposed to stamp each row with the date the
# forecast actually lands/publishes - the 15th of the month *after* the
# quarter start. Recompute every date in column A accordingly.
#
# Date math is done with plain Excel-serial arithmetic (Howard Hinnant's
# days_from_civil / civil_from_days, O(1), no calendar libraries needed)
# since this COM host's DateTime object support is limited.

function Convert-ExcelSerialToYmd($serial) {
    $z = [math]::Floor([double]$serial - 25569.0) + 719468
    if ($z -ge 0) {
        $era = [math]::Floor($z / 146097.0)
    } else {
        $era = [math]::Floor(($z - 146096.0) / 146097.0)
    }
    $doe = $z - $era * 146097
    $yoe = [math]::Floor(($doe - [math]::Floor($doe/1460) + [math]::Floor($doe/36524) - [math]::Floor($doe/146096)) / 365)
    $y = $yoe + $era * 400
    $doy = $doe - (365*$yoe + [math]::Floor($yoe/4) - [math]::Floor($yoe/100))
    $mp = [math]::Floor((5*$doy + 2) / 153)
    $d = $doy - [math]::Floor((153*$mp + 2)/5) + 1
    if ($mp -lt 10) { $m = $mp + 3 } else { $m = $mp - 9 }
    if ($m -le 2) { $y = $y + 1 }
    return @([int]$y, [int]$m, [int]$d)
}

function Convert-YmdToExcelSerial($y, $m, $d) {
    if ($m -le 2) { $yy = $y - 1 } else { $yy = $y }
    if ($yy -ge 0) {
        $era = [math]::Floor($yy / 400.0)
    } else {
        $era = [math]::Floor(($yy - 399) / 400.0)
    }
    $yoe = $yy - $era * 400
    if ($m -gt 2) { $mAdj = $m - 3 } else { $mAdj = $m + 9 }
    $doy = [math]::Floor((153*$mAdj + 2)/5) + $d - 1
    $doe = $yoe*365 + [math]::Floor($yoe/4) - [math]::Floor($yoe/100) + $doy
    $days0000 = $era*146097 + $doe
    $days1970 = $days0000 - 719468
    return [int]($days1970 + 25569)
}

function Get-RepublishSerial($serial) {
    $ymd = Convert-ExcelSerialToYmd $serial
    $y = $ymd[0]; $m = $ymd[1]
    $nm = $m + 1
    $ny = $y
    if ($nm -gt 12) {
        $nm = $nm - 12
        $ny = $ny + 1
    }
    return Convert-YmdToExcelSerial $ny $nm 15
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value2
    if ($old -ne $null) {
        $cell.Value = Get-RepublishSerial $old
    }
}
